# Generate Report for handback
#
# The 05daeeb2-... file has now been handed back (in sync with en-US),
# so its data row moves above the 097298bd-... row (which was already
# handed back earlier) on every sheet, and its "latest handback"
# timestamp is refreshed to reflect the just-completed handback.

$wb = $excel.ActiveWorkbook

# ====================== Sheet "Overview" ======================
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "05daeeb2-ac12-4594-84bc-3e5a63870673.md"
$ov.Range("B2").Value = "Handed back: in sync with en-US"
$ov.Range("C2").Value = "Handed back: in sync with en-US"

$ov.Range("A3").Value = "097298bd-7ea0-4fd0-9d70-83728cc19d14.md"
$ov.Range("B3").Value = "Handed back: in sync with en-US"
$ov.Range("C3").Value = "Handed back: in sync with en-US"

# Re-point the hyperlinks so the displayed text follows the swapped rows
# while keeping the very same relationship ids (and therefore the same
# underlying link targets) that existed before, in the same left-to-right,
# top-to-bottom order they were originally defined.
$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/af756787332082ec5e6400dc07cae1b2cb88645d/e2e/097298bd-7ea0-4fd0-9d70-83728cc19d14.md", "", "", "05daeeb2-ac12-4594-84bc-3e5a63870673.md")
$ov.Hyperlinks.Add($ov.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/b507373a6101f9fc256c88b1ce7ea876d0b91af5/e2e/05daeeb2-ac12-4594-84bc-3e5a63870673.md", "", "", "097298bd-7ea0-4fd0-9d70-83728cc19d14.md")
$ov.Hyperlinks.Add($ov.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/b507373a6101f9fc256c88b1ce7ea876d0b91af5/.localization-config", "", "", ".localization-config")

# ====================== Sheet "zh-cn" ======================
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "05daeeb2-ac12-4594-84bc-3e5a63870673.md"
$zh.Range("B2").Value = "Handed back: in sync with en-US"
$zh.Range("C2").Value = "05daeeb2-ac12-4594-84bc-3e5a63870673.25c69d1c765d57f5d8d6eccac1dd13afc922757f.zh-cn.xlf"
$zh.Range("D2").Value = "2016-01-14 03:15:02"
$zh.Range("E2").Value = "05daeeb2-ac12-4594-84bc-3e5a63870673.md"
$zh.Range("F2").Value = "05daeeb2-ac12-4594-84bc-3e5a63870673.25c69d1c765d57f5d8d6eccac1dd13afc922757f.zh-cn.xlf"
$zh.Range("G2").Value = "2016-01-14 03:15:50"
$zh.Range("H2").Value = "Include"

$zh.Range("A3").Value = "097298bd-7ea0-4fd0-9d70-83728cc19d14.md"
$zh.Range("B3").Value = "Handed back: in sync with en-US"
$zh.Range("C3").Value = "097298bd-7ea0-4fd0-9d70-83728cc19d14.f32b93eeefa6fc8adad04f453fa5c24b491312ad.zh-cn.xlf"
$zh.Range("D3").Value = "2016-01-14 03:13:15"
$zh.Range("E3").Value = "097298bd-7ea0-4fd0-9d70-83728cc19d14.md"
$zh.Range("F3").Value = "097298bd-7ea0-4fd0-9d70-83728cc19d14.f32b93eeefa6fc8adad04f453fa5c24b491312ad.zh-cn.xlf"
$zh.Range("G3").Value = "2016-01-14 03:13:59"
$zh.Range("H3").Value = "Include"

$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/af756787332082ec5e6400dc07cae1b2cb88645d/e2e/097298bd-7ea0-4fd0-9d70-83728cc19d14.md", "", "", "05daeeb2-ac12-4594-84bc-3e5a63870673.md")
$zh.Hyperlinks.Add($zh.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/745a7c469fbbcbc0aa50969d48f5ba3e387b4693/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/097298bd-7ea0-4fd0-9d70-83728cc19d14.f32b93eeefa6fc8adad04f453fa5c24b491312ad.zh-cn.xlf", "", "", "05daeeb2-ac12-4594-84bc-3e5a63870673.25c69d1c765d57f5d8d6eccac1dd13afc922757f.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/7bb567e4436655910eff5fa9e8653c211d152216/e2e/097298bd-7ea0-4fd0-9d70-83728cc19d14.md", "", "", "05daeeb2-ac12-4594-84bc-3e5a63870673.md")
$zh.Hyperlinks.Add($zh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/2611eafd55c728790e3bf2252d605e3a79d07197/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/097298bd-7ea0-4fd0-9d70-83728cc19d14.f32b93eeefa6fc8adad04f453fa5c24b491312ad.zh-cn.xlf", "", "", "05daeeb2-ac12-4594-84bc-3e5a63870673.25c69d1c765d57f5d8d6eccac1dd13afc922757f.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/b507373a6101f9fc256c88b1ce7ea876d0b91af5/e2e/05daeeb2-ac12-4594-84bc-3e5a63870673.md", "", "", "097298bd-7ea0-4fd0-9d70-83728cc19d14.md")
$zh.Hyperlinks.Add($zh.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8f2c89afa37d26fcbee6a2ba0b590fff2233bbd4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/05daeeb2-ac12-4594-84bc-3e5a63870673.25c69d1c765d57f5d8d6eccac1dd13afc922757f.zh-cn.xlf", "", "", "097298bd-7ea0-4fd0-9d70-83728cc19d14.f32b93eeefa6fc8adad04f453fa5c24b491312ad.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/7bb567e4436655910eff5fa9e8653c211d152216/e2e/05daeeb2-ac12-4594-84bc-3e5a63870673.md", "", "", "097298bd-7ea0-4fd0-9d70-83728cc19d14.md")
$zh.Hyperlinks.Add($zh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/2611eafd55c728790e3bf2252d605e3a79d07197/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/05daeeb2-ac12-4594-84bc-3e5a63870673.25c69d1c765d57f5d8d6eccac1dd13afc922757f.zh-cn.xlf", "", "", "097298bd-7ea0-4fd0-9d70-83728cc19d14.f32b93eeefa6fc8adad04f453fa5c24b491312ad.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/b507373a6101f9fc256c88b1ce7ea876d0b91af5/.localization-config", "", "", ".localization-config")

# ====================== Sheet "de-de" ======================
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "05daeeb2-ac12-4594-84bc-3e5a63870673.md"
$de.Range("B2").Value = "Handed back: in sync with en-US"
$de.Range("C2").Value = "05daeeb2-ac12-4594-84bc-3e5a63870673.25c69d1c765d57f5d8d6eccac1dd13afc922757f.de-de.xlf"
$de.Range("D2").Value = "2016-01-14 03:15:14"
$de.Range("E2").Value = "05daeeb2-ac12-4594-84bc-3e5a63870673.md"
$de.Range("F2").Value = "05daeeb2-ac12-4594-84bc-3e5a63870673.25c69d1c765d57f5d8d6eccac1dd13afc922757f.de-de.xlf"
$de.Range("G2").Value = "2016-01-14 03:16:11"
$de.Range("H2").Value = "Include"

$de.Range("A3").Value = "097298bd-7ea0-4fd0-9d70-83728cc19d14.md"
$de.Range("B3").Value = "Handed back: in sync with en-US"
$de.Range("C3").Value = "097298bd-7ea0-4fd0-9d70-83728cc19d14.f32b93eeefa6fc8adad04f453fa5c24b491312ad.de-de.xlf"
$de.Range("D3").Value = "2016-01-14 03:13:27"
$de.Range("E3").Value = "097298bd-7ea0-4fd0-9d70-83728cc19d14.md"
$de.Range("F3").Value = "097298bd-7ea0-4fd0-9d70-83728cc19d14.f32b93eeefa6fc8adad04f453fa5c24b491312ad.de-de.xlf"
$de.Range("G3").Value = "2016-01-14 03:14:21"
$de.Range("H3").Value = "Include"

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/af756787332082ec5e6400dc07cae1b2cb88645d/e2e/097298bd-7ea0-4fd0-9d70-83728cc19d14.md", "", "", "05daeeb2-ac12-4594-84bc-3e5a63870673.md")
$de.Hyperlinks.Add($de.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4675af6c2255ee046c2bb2d20562c04aa51a6779/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/097298bd-7ea0-4fd0-9d70-83728cc19d14.f32b93eeefa6fc8adad04f453fa5c24b491312ad.de-de.xlf", "", "", "05daeeb2-ac12-4594-84bc-3e5a63870673.25c69d1c765d57f5d8d6eccac1dd13afc922757f.de-de.xlf")
$de.Hyperlinks.Add($de.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/114fb4f6b350718aefcde4953525d9988a2ad73b/e2e/097298bd-7ea0-4fd0-9d70-83728cc19d14.md", "", "", "05daeeb2-ac12-4594-84bc-3e5a63870673.md")
$de.Hyperlinks.Add($de.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/d8c4707ab460a86380de43b407b29d5ed208ef1b/ol-handback/OpenLocalizationTestOrg/oltest.de-de/yuwzho/097298bd-7ea0-4fd0-9d70-83728cc19d14.f32b93eeefa6fc8adad04f453fa5c24b491312ad.de-de.xlf", "", "", "05daeeb2-ac12-4594-84bc-3e5a63870673.25c69d1c765d57f5d8d6eccac1dd13afc922757f.de-de.xlf")
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/b507373a6101f9fc256c88b1ce7ea876d0b91af5/e2e/05daeeb2-ac12-4594-84bc-3e5a63870673.md", "", "", "097298bd-7ea0-4fd0-9d70-83728cc19d14.md")
$de.Hyperlinks.Add($de.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fb4976e3f62c9ff00ed333f6dbc2d63dd696f792/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/05daeeb2-ac12-4594-84bc-3e5a63870673.25c69d1c765d57f5d8d6eccac1dd13afc922757f.de-de.xlf", "", "", "097298bd-7ea0-4fd0-9d70-83728cc19d14.f32b93eeefa6fc8adad04f453fa5c24b491312ad.de-de.xlf")
$de.Hyperlinks.Add($de.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/114fb4f6b350718aefcde4953525d9988a2ad73b/e2e/05daeeb2-ac12-4594-84bc-3e5a63870673.md", "", "", "097298bd-7ea0-4fd0-9d70-83728cc19d14.md")
$de.Hyperlinks.Add($de.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/d8c4707ab460a86380de43b407b29d5ed208ef1b/ol-handback/OpenLocalizationTestOrg/oltest.de-de/yuwzho/05daeeb2-ac12-4594-84bc-3e5a63870673.25c69d1c765d57f5d8d6eccac1dd13afc922757f.de-de.xlf", "", "", "097298bd-7ea0-4fd0-9d70-83728cc19d14.f32b93eeefa6fc8adad04f453fa5c24b491312ad.de-de.xlf")
$de.Hyperlinks.Add($de.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/b507373a6101f9fc256c88b1ce7ea876d0b91af5/.localization-config", "", "", ".localization-config")
